$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.686.02'
$ws.Range('E2').Value = '  +1.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.565.80'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.98'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.488'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '25.01'
$ws.Range('E8').Value = '  +5.87%  '
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0585'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0894'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.790.58'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.571.76'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.685.50'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.36'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.45'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.00'
$ws.Range('E23').Value = '  +0.78%  '
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.71'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.77'
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.23'
$ws.Range('E29').Value = '  -1.47%  '
$ws.Range('E30').Value = '  -3.95%  '
$ws.Range('E31').Value = '  -2.80%  '
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.394.34'
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('E34').Value = '  -3.00%  '
$ws.Range('E35').Value = '  -3.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.47'
$ws.Range('E37').Value = '  +1.88%  '
$ws.Range('E38').Value = '  -2.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.96'
$ws.Range('E40').Value = '  +2.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.519'
$ws.Range('E41').Value = '  -0.46%  '
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.771'
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0458'
$ws.Range('E44').Value = '  -3.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.92'
$ws.Range('E45').Value = '  +2.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.23'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.702.50'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.869'
$ws.Range('E48').Value = '  -5.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.12'
$ws.Range('E49').Value = '  -0.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.95'
$ws.Range('E50').Value = '  +5.47%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0100'
$ws.Range('E51').Value = '  -0.67%  '
